$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.029.18'
$ws.Range('E2').Value = '  -0.61%  '

$ws.Range('D3').Value = '1.619.17'
$ws.Range('E3').Value = '  -1.47%  '

$ws.Range('E4').Value = '  -0.10%  '

$ws.Range('D5').Value = '''214.43'
$ws.Range('E5').Value = '  -1.29%  '

$ws.Range('D6').Value = '''0.517'
$ws.Range('E6').Value = '  +0.52%  '

$ws.Range('E7').Value = '  -0.06%  '

$ws.Range('E8').Value = '  -1.50%  '

$ws.Range('D9').Value = '''0.0623'
$ws.Range('E9').Value = '  -0.56%  '

$ws.Range('D10').Value = '''20.24'
$ws.Range('E10').Value = '  +1.60%  '

$ws.Range('E11').Value = '  -0.37%  '

$ws.Range('D12').Value = '1.622.08'
$ws.Range('E12').Value = '  -1.84%  '

$ws.Range('E13').Value = '  -0.66%  '

$ws.Range('D14').Value = '''0.541'
$ws.Range('E14').Value = '  -0.47%  '

$ws.Range('D15').Value = '27.005.16'

$ws.Range('E16').Value = '  -4.53%  '

$ws.Range('D17').Value = '0.0₃0741'
$ws.Range('E17').Value = '  +0.39%  '

$ws.Range('D18').Value = '''215.83'
$ws.Range('E18').Value = '  -1.48%  '

$ws.Range('E19').Value = '  -0.07%  '

$ws.Range('D20').Value = '''6.88'
$ws.Range('E20').Value = '  +0.86%  '

$ws.Range('E21').Value = '  -0.76%  '

$ws.Range('D22').Value = '''2.41'
$ws.Range('E22').Value = '  -5.48%  '

$ws.Range('D23').Value = '''9.01'
$ws.Range('E23').Value = '  -1.80%  '

$ws.Range('D24').Value = '''147.19'
$ws.Range('E24').Value = '  -0.38%  '

$ws.Range('E25').Value = '  -0.12%  '

$ws.Range('D26').Value = '''7.27'
$ws.Range('E26').Value = '  -3.72%  '

$ws.Range('E27').Value = '  -0.92%  '

$ws.Range('D28').Value = '''15.54'
$ws.Range('E28').Value = '  -1.41%  '

$ws.Range('D29').Value = '''0.0503'
$ws.Range('E29').Value = '  -1.20%  '

$ws.Range('E30').Value = '  -1.10%  '

$ws.Range('D31').Value = '''3.35'

$ws.Range('E32').Value = '  -1.58%  '

$ws.Range('D33').Value = '1.335.30'
$ws.Range('E33').Value = '  +5.44%  '

$ws.Range('E34').Value = '  -0.70%  '

$ws.Range('E35').Value = '  -0.25%  '

$ws.Range('E36').Value = '  -1.14%  '

$ws.Range('D37').Value = '''0.544'
$ws.Range('E37').Value = '  -0.02%  '

$ws.Range('D38').Value = '''0.846'
$ws.Range('E38').Value = '  -0.14%  '

$ws.Range('E39').Value = '  -0.06%  '

$ws.Range('D40').Value = '''0.800'
$ws.Range('E40').Value = '  -1.05%  '

$ws.Range('E41').Value = '  +0.18%  '

$ws.Range('D42').Value = '''64.71'
$ws.Range('E42').Value = '  +4.89%  '

$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = '''5.21'
$ws.Range('E43').Value = '  -2.67%  '

$ws.Range('B44').Value = 'RocketPoolETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D44').Value = '1.756.51'
$ws.Range('E44').Value = '  -1.58%  '

$ws.Range('D45').Value = '''90.30'
$ws.Range('E45').Value = '  -1.49%  '

$ws.Range('E46').Value = '  +0.24%  '

$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').Value = '''0.838'
$ws.Range('E47').Value = '  +25.13%  '

$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0106'
$ws.Range('E48').Value = '  +3.69%  '

$ws.Range('E49').Value = '  -0.37%  '

$ws.Range('D50').Value = '''0.0996'
$ws.Range('E50').Value = '  +2.23%  '

$ws.Range('D51').Value = '''7.56'
$ws.Range('E51').Value = '  -0.68%  '
